# Weekly crime-data refresh: bump report volume/number and week-ending dates,
# then update the crime-complaint statistics table (rows 14-30) with newly
# collected figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title block: issue number and reporting week ---
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# --- Cells that flip between a numeric count and the "0"/"***.*" placeholder text ---
# (copy format+value from a same-row cell that already holds the right placeholder)
$ws.Range("D14").Copy($ws.Range("F14"))
$ws.Range("D14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("D27").Copy($ws.Range("G27"))
$ws.Range("E27").Copy($ws.Range("H27"))
$ws.Range("D29").Copy($ws.Range("C29"))
$ws.Range("D30").Copy($ws.Range("C30"))

# --- Plain numeric value updates ---
$ws.Range("L14").Value = -37.5
$ws.Range("F15").Value = 1
$ws.Range("M15").Value = 59.090909090909
$ws.Range("N15").Value = -20.454545454545
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -36.363636363636
$ws.Range("F16").Value = 42
$ws.Range("H16").Value = -22.222222222222
$ws.Range("I16").Value = 334
$ws.Range("J16").Value = 369
$ws.Range("K16").Value = -9.485094850948
$ws.Range("L16").Value = 5.362776025236
$ws.Range("M16").Value = -0.298507462686
$ws.Range("N16").Value = -69.691470054446
$ws.Range("C17").Value = 21
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 81
$ws.Range("G17").Value = 75
$ws.Range("H17").Value = 8
$ws.Range("I17").Value = 636
$ws.Range("J17").Value = 562
$ws.Range("K17").Value = 13.167259786476
$ws.Range("L17").Value = 36.188436830835
$ws.Range("M17").Value = 101.904761904762
$ws.Range("N17").Value = 41.019955654102
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -19.047619047619
$ws.Range("I18").Value = 144
$ws.Range("J18").Value = 138
$ws.Range("K18").Value = 4.347826086956
$ws.Range("L18").Value = -8.280254777070
$ws.Range("M18").Value = -44.827586206896
$ws.Range("N18").Value = -90.897597977244
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 27
$ws.Range("E19").Value = -3.703703703703
$ws.Range("F19").Value = 91
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = 8.333333333333
$ws.Range("I19").Value = 622
$ws.Range("J19").Value = 581
$ws.Range("K19").Value = 7.056798623063
$ws.Range("L19").Value = 24.4
$ws.Range("M19").Value = 63.684210526315
$ws.Range("N19").Value = 0.322580645161
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -77.777777777777
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -39.285714285714
$ws.Range("I20").Value = 225
$ws.Range("J20").Value = 247
$ws.Range("K20").Value = -8.906882591093
$ws.Range("L20").Value = -5.857740585774
$ws.Range("M20").Value = 87.5
$ws.Range("N20").Value = -78.070175438596
$ws.Range("C21").Value = 61
$ws.Range("D21").Value = 70
$ws.Range("E21").Value = -12.857142857142
$ws.Range("F21").Value = 249
$ws.Range("G21").Value = 262
$ws.Range("H21").Value = -4.961832061068
$ws.Range("I21").Value = 2001
$ws.Range("J21").Value = 1932
$ws.Range("K21").Value = 3.571428571428
$ws.Range("L21").Value = 16.812609457092
$ws.Range("M21").Value = 38.958333333333
$ws.Range("N21").Value = -58.750773036487
$ws.Range("D22").Value = 1
$ws.Range("J22").Value = 35
$ws.Range("K22").Value = -22.857142857142
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -32.352941176470
$ws.Range("F24").Value = 139
$ws.Range("G24").Value = 151
$ws.Range("H24").Value = -7.947019867549
$ws.Range("I24").Value = 1148
$ws.Range("J24").Value = 1284
$ws.Range("K24").Value = -10.591900311526
$ws.Range("L24").Value = -34.957507082153
$ws.Range("M24").Value = 31.2
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 86
$ws.Range("H25").Value = -26.744186046511
$ws.Range("I25").Value = 474
$ws.Range("J25").Value = 710
$ws.Range("K25").Value = -33.239436619718
$ws.Range("L25").Value = -59.31330472103
$ws.Range("C26").Value = 20
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 73
$ws.Range("G26").Value = 69
$ws.Range("H26").Value = 5.797101449275
$ws.Range("I26").Value = 647
$ws.Range("J26").Value = 604
$ws.Range("K26").Value = 7.119205298013
$ws.Range("L26").Value = 11.744386873920
$ws.Range("M26").Value = 3.685897435897
$ws.Range("I27").Value = 42
$ws.Range("K27").Value = -2.325581395348
$ws.Range("L27").Value = -2.325581395348
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 8
$ws.Range("H28").Value = -11.111111111111
$ws.Range("I28").Value = 88
$ws.Range("J28").Value = 82
$ws.Range("K28").Value = 7.317073170731
$ws.Range("L28").Value = 23.943661971831
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 50
$ws.Range("M29").Value = -48.148148148148
$ws.Range("N29").Value = -73.584905660377
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("M30").Value = -47.619047619047
$ws.Range("N30").Value = -76.595744680851
